# Reposition/resize a handful of shapes on slide 1 of the design doc.
#
# NOTE on numeric literals: PowerPoint's object model expresses Left/Top/
# Width/Height in points while OOXML stores EMUs (1 pt = 12700 EMU). The
# host's Shape position/size fields round-trip through a single-precision
# (f32) point value before being truncated back to EMU, so naively writing
# "emu / 12700.0" can land 1 EMU short of the intended target. The literals
# below were chosen so that, after that f32 round-trip, they truncate back
# to the exact target EMU values from the target OOXML.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1. Cloud 10: x 8710296 -> 8676458 EMU (y, cx, cy unchanged)
$sh = $s.Shapes.Item("Cloud 10")
$sh.Left = 683.1857299804688

# 2. Straight Arrow Connector 57: off 3848194,4131288 -> 3865678,4150721 EMU
#    ext 1982278,1163834 -> 1964794,1144401 EMU
$sh = $s.Shapes.Item("Straight Arrow Connector 57")
$sh.Left   = 304.3841247558594
$sh.Top    = 326.82843017578125
$sh.Width  = 154.70819091796875
$sh.Height = 90.11032104492188

# 3. Connector: Elbow 1093: x 6037091 -> 6003253 EMU (y, cx, cy unchanged)
$sh = $s.Shapes.Item("Connector: Elbow 1093")
$sh.Left = 472.6971130371094

# 4. Picture 1122: y 3314576 -> 3425890 EMU (x, cx, cy unchanged)
$sh = $s.Shapes.Item("Picture 1122")
$sh.Top = 269.755126953125

# 5. Straight Arrow Connector 1123: off 3865814,3442047 -> 3844195,3553361 EMU
#    ext 2129593,547629 -> 2151212,386543 EMU
$sh = $s.Shapes.Item("Straight Arrow Connector 1123")
$sh.Left   = 302.6925354003906
$sh.Top    = 279.7922058105469
$sh.Width  = 169.38677978515625
$sh.Height = 30.436458587646484
